$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-10-04 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-05 Thursday", 2) | Out-Null
$d.Content.Find.Execute("24×56=1344", $true, $false, $false, $false, $false, $true, 1, $false, "11×12=132", 2) | Out-Null
$d.Content.Find.Execute("94×80=7520", $true, $false, $false, $false, $false, $true, 1, $false, "59×98=5782", 2) | Out-Null
$d.Content.Find.Execute("21×49=1029", $true, $false, $false, $false, $false, $true, 1, $false, "56×34=1904", 2) | Out-Null
$d.Content.Find.Execute("52×19=988", $true, $false, $false, $false, $false, $true, 1, $false, "99×85=8415", 2) | Out-Null
$d.Content.Find.Execute("76×19=1444", $true, $false, $false, $false, $false, $true, 1, $false, "92×44=4048", 2) | Out-Null
$d.Content.Find.Execute("69×56=3864", $true, $false, $false, $false, $false, $true, 1, $false, "53×78=4134", 2) | Out-Null
$d.Content.Find.Execute("59×39=2301", $true, $false, $false, $false, $false, $true, 1, $false, "65×68=4420", 2) | Out-Null
$d.Content.Find.Execute("70×39=2730", $true, $false, $false, $false, $false, $true, 1, $false, "66×86=5676", 2) | Out-Null
$d.Content.Find.Execute("93×44=4092", $true, $false, $false, $false, $false, $true, 1, $false, "17×47=799", 2) | Out-Null
$d.Content.Find.Execute("70×21=1470", $true, $false, $false, $false, $false, $true, 1, $false, "38×21=798", 2) | Out-Null
$d.Content.Find.Execute("13×47=611", $true, $false, $false, $false, $false, $true, 1, $false, "38×12=456", 2) | Out-Null
$d.Content.Find.Execute("79×69=5451", $true, $false, $false, $false, $false, $true, 1, $false, "31×64=1984", 2) | Out-Null
$d.Content.Find.Execute("98×60=5880", $true, $false, $false, $false, $false, $true, 1, $false, "30×91=2730", 2) | Out-Null
$d.Content.Find.Execute("53×29=1537", $true, $false, $false, $false, $false, $true, 1, $false, "84×32=2688", 2) | Out-Null
$d.Content.Find.Execute("27×97=2619", $true, $false, $false, $false, $false, $true, 1, $false, "28×11=308", 2) | Out-Null
$d.Content.Find.Execute("39×47=1833", $true, $false, $false, $false, $false, $true, 1, $false, "16×38=608", 2) | Out-Null
$d.Content.Find.Execute("76×14=1064", $true, $false, $false, $false, $false, $true, 1, $false, "18×90=1620", 2) | Out-Null
$d.Content.Find.Execute("37×21=777", $true, $false, $false, $false, $false, $true, 1, $false, "25×16=400", 2) | Out-Null
$d.Content.Find.Execute("42×29=1218", $true, $false, $false, $false, $false, $true, 1, $false, "85×40=3400", 2) | Out-Null
$d.Content.Find.Execute("66×76=5016", $true, $false, $false, $false, $false, $true, 1, $false, "37×11=407", 2) | Out-Null
$d.Content.Find.Execute("94×64=6016", $true, $false, $false, $false, $false, $true, 1, $false, "34×46=1564", 2) | Out-Null
$d.Content.Find.Execute("43×20=860", $true, $false, $false, $false, $false, $true, 1, $false, "35×47=1645", 2) | Out-Null
$d.Content.Find.Execute("30×25=750", $true, $false, $false, $false, $false, $true, 1, $false, "42×97=4074", 2) | Out-Null
$d.Content.Find.Execute("46×21=966", $true, $false, $false, $false, $false, $true, 1, $false, "11×34=374", 2) | Out-Null
$d.Content.Find.Execute("32×43=1376", $true, $false, $false, $false, $false, $true, 1, $false, "62×74=4588", 2) | Out-Null
